{"js": "// Fix duplicated \"vistos\" text in the 2022 resolution addendum header template.\n//\n// Original (buggy) wording had the placeholder `${art8}` immediately followed\n// by a duplicated sentence (\"Decreto N\u00b0140/04, del Ministerio de Salud que\n// aprob\u00f3 el Reglamento Org\u00e1nico de los Servicios de Salud;\") before the\n// \"Ley N\u00b019.880 ...\" clause. The fix removes that duplicated sentence\n// (together with the now orphaned `${art8}` placeholder) and re-inserts\n// `${art8}` \u2014 with its original yellow highlight \u2014 right before\n// `${directorDecreto}`, i.e. at the end of the \"vistos\" citation list.\n\nconst body = context.document.body;\n\n// 1) Locate and delete \"${art8} Decreto N\u00b0140/04, ... de los Servicios de Salud;\"\nconst dupSearch = body.search(\n  \"${art8} Decreto N\u00b0140/04, del Ministerio de Salud que aprob\u00f3 el Reglamento Org\u00e1nico de los Servicios de Salud;\",\n  { matchCase: true }\n);\ndupSearch.load(\"items\");\nawait context.sync();\n\nif (dupSearch.items.length > 0) {\n  dupSearch.items[0].delete();\n  await context.sync();\n}\n\n// 2) Re-insert \"${art8} \" (highlighted yellow) right before \"${directorDecreto}\"\nconst targetSearch = body.search(\"${directorDecreto}\", { matchCase: true });\ntargetSearch.load(\"items\");\nawait context.sync();\n\nif (targetSearch.items.length > 0) {\n  const target = targetSearch.items[0];\n  const inserted = target.insertText(\"${art8} \", Word.InsertLocation.before);\n  inserted.font.highlightColor = \"yellow\";\n  await context.sync();\n}\n", "ps1": "# Fix duplicated \"vistos\" text in the 2022 resolution addendum header template.\n#\n# Original (buggy) wording had the placeholder ${art8} immediately followed\n# by a duplicated sentence (\"Decreto N\u00b0140/04, del Ministerio de Salud que\n# aprob\u00f3 el Reglamento Org\u00e1nico de los Servicios de Salud;\") before the\n# \"Ley N\u00b019.880 ...\" clause. The fix removes that duplicated sentence\n# (together with the now orphaned ${art8} placeholder) and re-inserts\n# ${art8} - with its original yellow highlight - right before\n# ${directorDecreto}, i.e. at the end of the \"vistos\" citation list.\n#\n# NOTE: literal ${...} template placeholders must be written inside\n# *single*-quoted strings below, otherwise PowerShell parses them as\n# variable references (and substitutes an empty string).\n\n$d = $word.ActiveDocument\n\n# 1) Locate and delete \"${art8} Decreto N\u00b0140/04, ... de los Servicios de Salud;\"\n$dup = $d.Content\n$foundDup = $dup.Find.Execute('${art8} Decreto N\u00b0140/04, del Ministerio de Salud que aprob\u00f3 el Reglamento Org\u00e1nico de los Servicios de Salud;')\nif ($foundDup) {\n    $dup.Delete()\n}\n\n# 2) Re-insert \"${art8} \" right before \"${directorDecreto}\"\n$target = $d.Content\n$foundTarget = $target.Find.Execute('${directorDecreto}')\nif ($foundTarget) {\n    $target.Collapse(1)  # wdCollapseStart\n    $target.InsertBefore('${art8} ')\n}\n\n# 3) Highlight the re-inserted \"${art8}\" placeholder in yellow (use Font.\n#    HighlightColorIndex, scoped to the run's own font, rather than\n#    Range.HighlightColorIndex which bleeds into the rest of the story).\n$art8Range = $d.Content\n$foundArt8 = $art8Range.Find.Execute('${art8}')\nif ($foundArt8) {\n    $art8Range.Font.HighlightColorIndex = 7  # wdYellow\n}\n"}
